# Apply the "mise a jour des notebook residentiel_tertiaire" edit.
#
# Summary of the change:
#  - the worksheet "res_type_Energy_source_year" is removed entirely
#  - the worksheet "year_res_type" moves up to take its old slot
#  - a brand-new worksheet "year" is appended at the end, holding a
#    year -> retrofit_change_total_proportion_surface curve
#  - on sheet "0D" the row holding
#    "retrofit_change_total_proportion_surface" / 1 is deleted (that
#    hypothesis is now represented by the new "year" sheet instead of a
#    single flat number)
#  - sheet "0D" becomes the active / selected tab instead of
#    "res_type_Energy_source"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Remove the obsolete "res_type_Energy_source_year" sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("res_type_Energy_source_year").Delete()

# ---------------------------------------------------------------------
# 2) Sheet "0D": drop the row for retrofit_change_total_proportion_surface
#    (Name column A, Value column B) - it becomes the "year" curve sheet.
# ---------------------------------------------------------------------
$ws0D = $wb.Worksheets.Item("0D")

$foundRow = 0
for ($r = 1; $r -le $ws0D.Cells.Item($ws0D.Rows.Count, 1).End(-4162).Row; $r++) {
    if ($ws0D.Cells.Item($r, 1).Value2 -eq "retrofit_change_total_proportion_surface") {
        $foundRow = $r
    }
}
if ($foundRow -gt 0) {
    $ws0D.Rows($foundRow).Delete()
}

# ---------------------------------------------------------------------
# 3) Add the new "year" worksheet at the very end of the workbook, right
#    after "year_res_type", with the year -> proportion curve.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$yearWs = $wb.Worksheets.Add($null, $lastSheet)
$yearWs.Name = "year"

$yearWs.Range("A1").Value = "year"
$yearWs.Range("B1").Value = "retrofit_change_total_proportion_surface"

$years = @(2020, 2022, 2025, 2030, 2035, 2040, 2045, 2050)
$vals  = @(0, 0.005, 0.1, 0.35, 0.7, 0.85, 0.95, 1)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $yearWs.Cells.Item($row, 1).Value = $years[$i]
    $yearWs.Cells.Item($row, 2).Value = $vals[$i]
}

# formatting: bold header + size-12 body font, thin borders under the
# header row and around the data column
$headerRange = $yearWs.Range("A1:B1")
$bodyRange   = $yearWs.Range("A2:B9")
$allRange    = $yearWs.Range("A1:B9")

$allRange.Font.Bold = $true
$allRange.Font.Size = 11
$headerRange.Font.Size = 11

$allRange.Borders.Item(7).LineStyle = 1
$allRange.Borders.Item(7).Weight = 2
$allRange.Borders.Item(10).LineStyle = 1
$allRange.Borders.Item(10).Weight = 2
$headerRange.Borders.Item(9).LineStyle = 1
$headerRange.Borders.Item(9).Weight = 2
for ($i = 2; $i -le 9; $i++) {
    $yearWs.Range("A$i`:B$i").Borders.Item(9).LineStyle = 1
    $yearWs.Range("A$i`:B$i").Borders.Item(9).Weight = 2
}

$yearWs.Columns.Item(1).ColumnWidth = 14.5
for ($r = 1; $r -le 9; $r++) {
    $yearWs.Rows($r).RowHeight = 16
}

# ---------------------------------------------------------------------
# 4) Make "0D" the active / selected sheet (it used to be
#    "res_type_Energy_source").
# ---------------------------------------------------------------------
$ws0D.Activate()
$ws0D.Range("A1").Select()
